$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; this shifts all existing rows 75-198 down to 76-199,
# preserving all of their data exactly as before (matching the diff's row-shift pattern).
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75 with the new record.
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44557
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 100112003
$ws.Range("G75").Value = "Ajo"
$ws.Range("H75").Value = "Chino"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 20
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 21000
$ws.Range("M75").Value = 20500
$ws.Range("N75").Value = "$/caja 10 kilos"
$ws.Range("O75").Value = "China"
$ws.Range("P75").Value = 2050
$ws.Range("Q75").Value = 10
$ws.Range("R75").Value = "Hortaliza"
